$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33.28590577994189
$ws.Range("C2").Value = 16.275184844542355
$ws.Range("D2").Value = 0.48895123816488695
$ws.Range("E2").Value = 30.43706377061801
$ws.Range("F2").Value = 15.354937458636416
$ws.Range("G2").Value = 0.5044815615052556
$ws.Range("H2").Value = 281.5
$ws.Range("I2").Value = 247.5
$ws.Range("B3").Value = 33.31607799695591
$ws.Range("C3").Value = 16.30342119480783
$ws.Range("D3").Value = 0.48935595589305175
$ws.Range("E3").Value = 30.501523533415934
$ws.Range("F3").Value = 15.415012223537902
$ws.Range("G3").Value = 0.5053849918890114
$ws.Range("H3").Value = 281.5
$ws.Range("I3").Value = 247.5
$ws.Range("B4").Value = 33.34111401719178
$ws.Range("C4").Value = 16.32864032274177
$ws.Range("D4").Value = 0.4897448931767002
$ws.Range("E4").Value = 30.49446724311285
$ws.Range("F4").Value = 15.397149174403538
$ws.Range("G4").Value = 0.5049161558276106
$ws.Range("H4").Value = 281.5
$ws.Range("I4").Value = 247.5
$ws.Range("B5").Value = 33.36158903980185
$ws.Range("C5").Value = 16.348184105031642
$ws.Range("D5").Value = 0.4900301387181389
$ws.Range("E5").Value = 30.515752826899742
$ws.Range("F5").Value = 15.399669414184094
$ws.Range("G5").Value = 0.5046465509646294
$ws.Range("H5").Value = 281.5
$ws.Range("I5").Value = 247.5
$ws.Range("B6").Value = 33.37600731643254
$ws.Range("C6").Value = 16.307877027489013
$ws.Range("D6").Value = 0.48861078177736067
$ws.Range("E6").Value = 30.54521312640602
$ws.Range("F6").Value = 15.415925728397664
$ws.Range("G6").Value = 0.5046920335635424
$ws.Range("H6").Value = 281.5
$ws.Range("I6").Value = 247
$ws.Range("B7").Value = 33.38970897133535
$ws.Range("C7").Value = 16.311471467215917
$ws.Range("D7").Value = 0.4885179287192683
$ws.Range("E7").Value = 30.587524235588717
$ws.Range("F7").Value = 15.482524661983135
$ws.Range("G7").Value = 0.5061712266326266
$ws.Range("H7").Value = 281.5
$ws.Range("I7").Value = 247
$ws.Range("B8").Value = 33.39952333340587
$ws.Range("C8").Value = 16.322743477919445
$ws.Range("D8").Value = 0.4887118691778934
$ws.Range("E8").Value = 30.602007772780688
$ws.Range("F8").Value = 15.494615353011225
$ws.Range("G8").Value = 0.5063267569911897
$ws.Range("H8").Value = 281.5
$ws.Range("I8").Value = 247
$ws.Range("B9").Value = 33.40974669814801
$ws.Range("C9").Value = 16.329580341143316
$ws.Range("D9").Value = 0.488766960392684
$ws.Range("E9").Value = 30.592300539390152
$ws.Range("F9").Value = 15.488203735858297
$ws.Range("G9").Value = 0.506277836670568
$ws.Range("H9").Value = 281.5
$ws.Range("I9").Value = 246.5
$ws.Range("B10").Value = 33.41907340895121
$ws.Range("C10").Value = 16.34288652076675
$ws.Range("D10").Value = 0.4890287148532775
$ws.Range("E10").Value = 30.615062913747778
$ws.Range("F10").Value = 15.476099470649883
$ws.Range("G10").Value = 0.5055060482564058
$ws.Range("H10").Value = 281.5
$ws.Range("I10").Value = 246.5
$ws.Range("B11").Value = 33.428031668027195
$ws.Range("C11").Value = 16.34890624067611
$ws.Range("D11").Value = 0.4890777417897835
$ws.Range("E11").Value = 30.61591614602488
$ws.Range("F11").Value = 15.457927257041952
$ws.Range("G11").Value = 0.504898405891701
$ws.Range("H11").Value = 281.5
$ws.Range("I11").Value = 246.5

$ws.Columns.Item(1).ColumnWidth = 38.5
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 31.0
$ws.Columns.Item(4).ColumnWidth = 37.833333333333336
$ws.Columns.Item(5).ColumnWidth = 31.166666666666668
$ws.Columns.Item(6).ColumnWidth = 29.666666666666668
$ws.Columns.Item(7).ColumnWidth = 36.666666666666664
$ws.Columns.Item(8).ColumnWidth = 31.5
$ws.Columns.Item(9).ColumnWidth = 30.166666666666668
